$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top (shifts existing rows down by 3).
$ws.Range("A1:A3").EntireRow.Insert()

# The insert does not carry the date number format into the new rows;
# copy it from the row that already has it (now row 4, the old row 1)
# so A1:A3 end up on the same style index as the rest of column A.
$ws.Range("A4").Copy()
$ws.Range("A1:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- New row 1: 2014-04-23 ---
$ws.Range("A1").Value2 = 41752
$ws.Range("B1").Value2 = "INTERES A SU FAVOR"
$ws.Range("C1").Value2 = "C"
$ws.Range("D1").Value2 = "0000858352"
$ws.Range("E1").Value2 = "AGENCIA PARA PROCESOS BATCH"
$ws.Range("F1").Value2 = "0.26  "
$ws.Range("G1").Value2 = "4140.36"

# --- New row 2: 2014-04-22 ---
$ws.Range("A2").Value2 = 41751
$ws.Range("B2").Value2 = "INTERES A SU FAVOR"
$ws.Range("C2").Value2 = "C"
$ws.Range("D2").Value2 = "0000858357"
$ws.Range("E2").Value2 = "AGENCIA PARA PROCESOS BATCH"
$ws.Range("F2").Value2 = "0.26  "
$ws.Range("G2").Value2 = "4140.10"

# --- New row 3: 2014-04-21 ---
$ws.Range("A3").Value2 = 41750
$ws.Range("B3").Value2 = "INTERES A SU FAVOR"
$ws.Range("C3").Value2 = "C"
$ws.Range("D3").Value2 = "0000858366"
$ws.Range("E3").Value2 = "AGENCIA PARA PROCESOS BATCH"
$ws.Range("F3").Value2 = "0.26  "
$ws.Range("G3").Value2 = "4139.84"

# Shared formula across the new rows (H1:H3): updated time format plus the
# two new array keys ('mo_fecha_borrado', 'mo_quien_borra').
$formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A1,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B1,"'', ''mo_tipo'' => ''",C1,"'', ''mo_documento'' => ''",D1,"'', ''mo_oficina'' => ''",E1,"'', ''mo_monto'' => ",F1,", ''mo_saldo'' => ",G1,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd H:m:s"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL, ''mo_fecha_borrado'' => NULL, ''mo_quien_borra'' => NULL, ''mo_borrado_logico'' => false),")'
$ws.Range("H1:H3").Formula = $formula

# The old shared formula used to span H1:H5; after the insert those cells
# are now H4:H8 and must be cleared since only H1:H3 keep the formula now.
$ws.Range("H4:H8").ClearContents()

# Match the updated selection shown in the workbook view.
$ws.Range("H1:H3").Select()
